# Update handback-status timestamps as part of "Generate Report for Handback".
$wb = $excel.ActiveWorkbook

# Sheet "Overview": Latest HO Xliff Generate Date (G2)
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G2").Value = "2016-08-27 09:04:58"

# Sheet "zh-cn": Correspond Handoff Datetime (H2), Correspond Handback DateTime (K2)
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H2").Value = "2016-08-27 09:04:54"
$wsZhCn.Range("K2").Value = "2016-08-27 09:05:15"

# Sheet "de-de": Correspond Handback DateTime (K2)
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("K2").Value = "2016-08-27 09:05:22"
